$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.124.65'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '1.836.11'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6282'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07524'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2926'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  +2.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07711'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '1.831.56'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.006'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6690'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009384'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.996'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '29.134.98'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").Value = '2.075.75'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '223.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.006'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.77%  '
$ws.Range("E23").Value = '  -0.93%  '
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1396'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.511'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.500'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05668'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.157'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.36%  '
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.205'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7473'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.847'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.765'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").Value = '1.221.81'
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01785'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  +2.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8941'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000126'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.976.92'
$ws.Range("E46").Value = '  -0.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.07680'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +12.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5096'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4083'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.003'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.48%  '
